$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The D3:D9 cells are formatted as Text (numFmtId 49, "@"), so writing a
# numeric literal to .Value while that format is active gets stored as a
# text string (matches real Excel UI behaviour). Temporarily switch the
# range to General, write the new numbers, then restore the Text format
# so the numbers are written as proper numeric cells (t="n"), matching
# the original file's cell typing.
$priceRange = $ws.Range("D3:D9")
$priceRange.NumberFormat = "General"

$ws.Range("D3").Value = 263.86
$ws.Range("D4").Value = 146750
$ws.Range("D5").Value = 22766
$ws.Range("D6").Value = 778.1
$ws.Range("D7").Value = 221.88
$ws.Range("D8").Value = 843.4
$ws.Range("D9").Value = 40.15

$priceRange.NumberFormat = "@"

# Remove the trailing empty rows 49 and 50 (shrinks the used range/dimension
# from A1:F50 down to A1:F48).
$ws.Rows("49:50").Delete()

# Update the active selection from B3:B48 to A3:A48.
$ws.Range("A3:A48").Select()
